$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..70 (row 1 is the header and is unchanged).
# The exercise list was re-synced: several new exercises/variants were
# inserted, a few labels were renamed/cleaned up, and two rows were
# removed from the old "Quads" section, which reflows every row below
# the first insertion point.
$data = @(
    @('Chest', 'Barbell Bench Press ', 'Overall Chest ', 'Barbell, Bench'),
    @('Chest', 'Incline Barbell Press', 'Upper Chest', 'Barbell, Incline Bench'),
    @('Chest', 'Decline Barbell Press', 'Lower Chest', 'Barbell, Decline Bench'),
    @('Chest', 'Dumbbell Chest Press', 'Overall Chest ', 'Dumbbells, Bench'),
    @('Chest', 'Incline Dumbbeell Chest Press', 'Upper Chest', 'Dumbells, Incline Bench'),
    @('Chest', 'Incline Chest Press Machine', 'Upper Chest', 'Machine'),
    @('Chest', 'Chest Press Machine', 'Overall Chest ', 'Machine'),
    @('Chest', 'Plate Loaded Chest Press', 'Overall Chest ', 'Machine'),
    @('Chest', 'Dips', 'Lower Chest / Triceps', 'Dip Station, Belt'),
    @('Chest', 'Cable Fly', 'Chest', 'Cable Machine'),
    @('Chest', 'High-to-Low Cable Fly', 'Chest', 'Cable Machine'),
    @('Chest', 'Low-to-High Cable Fly', 'Chest', 'Cable Machine'),
    @('Chest', 'Pec Deck', 'Chest', 'Machine'),
    @('Chest', 'Dumbbell Pullover', 'Serratus / Upper Chest', 'Dumbbell, Bench'),
    @('Chest', 'Landmine Press', 'Upper Chest', 'Landmine Attachment'),
    @('Back', 'Deadlift (Conventional)', 'Posterior Chain / Lats', 'Barbell'),
    @('Back', 'Pull-Up ', 'Lats ', 'Pull-Up Bar'),
    @('Back', 'Chin-Up', 'Lats / Biceps', 'Pull-Up Bar'),
    @('Back', 'Bent-Over Row (Barbell)', 'Back Thickness', 'Barbell'),
    @('Back', 'T-Bar Row', 'Mid-Back Thickness', 'T-Bar Station'),
    @('Back', 'One-Arm Dumbbell Row', 'Lats (Unilateral)', 'Dumbbell, Bench'),
    @('Back', 'Seated Cable Row', 'Mid-Back / Lats', 'Cable Machine'),
    @('Back', 'Cable Wide Row', 'Mid-Back / Lats', 'Cable Machine'),
    @('Back', 'Lat Pulldown ', 'Lats', 'Cable Machine'),
    @('Back', 'Straight-Arm Pulldown', 'Lats (Isolation)', 'Cable Machine'),
    @('Back', 'Face Pull', 'Rear Delts / Rotators', 'Cable Machine'),
    @('Back', 'Barbell Shrug', 'Traps', 'Barbell'),
    @('Back', 'Dumbbell Shrugs', 'Traps', 'Dumbbell'),
    @('Shoulders', 'Overhead Press (Military)', 'Overall Shoulder Mass', 'Barbell'),
    @('Shoulders', 'Seated Dumbbell Press', 'Front / Side Delts', 'Dumbbells, Bench'),
    @('Shoulders', 'Arnold Press', 'Front / Side Delts', 'Dumbbells'),
    @('Shoulders', 'Lateral Raise', 'Side Delts (Width)', 'Dumbbells / Cables'),
    @('Shoulders', 'Front Raise', 'Front Delts', 'Dumbbells / Plate'),
    @('Shoulders', 'Rear Delt Flye (Reverse)', 'Rear Delts', 'Machine / Dumbbells'),
    @('Shoulders', 'Upright Row', 'Side Delts / Traps', 'Barbell / Cables'),
    @('Quads', 'Barbell Squat', 'Overall Leg Mass', 'Barbell, Rack'),
    @('Quads', 'Front Squat', 'Quads / Core', 'Barbell, Rack'),
    @('Quads', 'Leg Press', 'Quads ', 'Machine'),
    @('Quads', 'Hack Squat', 'Quads ', 'Machine'),
    @('Quads', 'Bulgarian Split Squat', 'Quads / Glutes', 'Dumbbells / Barbell'),
    @('Quads', 'Lunges (Walking/Reverse)', 'Quads / Glutes', 'Dumbbells / Barbell'),
    @('Quads', 'Leg Extension', 'Quads (Isolation)', 'Machine'),
    @('Hams/Glutes', 'Romanian Deadlift (RDL)', 'Hamstrings / Glutes', 'Barbell / Dumbbells'),
    @('Hams/Glutes', 'Stiff-Legged Deadlift', 'Hamstrings (Stretch)', 'Barbell'),
    @('Hams/Glutes', 'Hip Thrust', 'Glutes (Max Load)', 'Barbell, Bench'),
    @('Hams/Glutes', 'Glute Bridge', 'Glutes', 'Barbell / Bodyweight'),
    @('Hams/Glutes', 'Lying Leg Curl', 'Hamstrings', 'Machine'),
    @('Hams/Glutes', 'Seated Leg Curl', 'Hamstrings', 'Machine'),
    @('Hams/Glutes', 'Glute-Ham Raise (GHR)', 'Posterior Chain', 'GHR Bench'),
    @('Hams/Glutes', 'Cable Pull-Through', 'Glutes / Hamstrings', 'Cable Machine'),
    @('Calves', 'Standing Calf Raise', 'Gastrocnemius', 'Machine / Smith'),
    @('Calves', 'Seated Calf Raise', 'Soleus', 'Machine'),
    @('Triceps', 'Close-Grip Bench Press', 'Overall Triceps', 'Barbell, Bench'),
    @('Triceps', 'Skullcrushers', 'Long Head', 'EZ-Bar / Barbell'),
    @('Triceps', 'Tricep Pushdown ', 'Lateral Head', 'Cable Machine'),
    @('Triceps', 'Overhead Extension', 'Long Head (Stretch)', 'Dumbbell / Cable'),
    @('Triceps', 'Tricep Kickback', 'Contraction', 'Dumbbell / Cable'),
    @('Biceps', 'Barbell Curl', 'Overall Mass', 'Barbell'),
    @('Biceps', 'Dumbbell Curl ', 'Biceps', 'Dumbbells'),
    @('Biceps', 'Hammer Curl', 'Brachialis / Forearm', 'Dumbbells'),
    @('Biceps', 'Preacher Curl', 'Short Head', 'EZ-Bar / Machine'),
    @('Biceps', 'Incline Dumbbell Curl', 'Long Head (Stretch)', 'Dumbbells, Incline Bench'),
    @('Biceps', 'Concentration Curl', 'Peak / Isolation', 'Dumbbell'),
    @('Abs/Core', 'Hanging Leg Raise', 'Abs', 'Pull-Up Bar'),
    @('Abs/Core', 'Cable Crunch', 'Abs', 'Cable Machine'),
    @('Abs/Core', 'Decline Crunch', 'Abs', 'Decline Bench'),
    @('Abs/Core', 'Ab Wheel Rollout', 'Core Stability', 'Ab Wheel'),
    @('Abs/Core', 'Woodchoppers', 'Obliques', 'Cable Machine'),
    @('Abs/Core', 'Plank (Weighted)', 'Isometric Stability', 'Plate / Bodyweight')
)

$firstDataRow = 2
$lastDataRow = $firstDataRow + $data.Length - 1   # 70
$existingLastRow = 62

# Make sure any newly-added rows (63-70) inherit the same formatting
# (font/style) as the rest of the table before we write values into
# them, by copying the format of an already-styled row.
if ($lastDataRow -gt $existingLastRow) {
    $ws.Range("A2:D2").Copy() | Out-Null
    $newRowsRange = $ws.Range("A" + ($existingLastRow + 1) + ":D" + $lastDataRow)
    $newRowsRange.PasteSpecial(-4122) | Out-Null
}

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $firstDataRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowValues[$j]
    }
}

Write-Host "Updated $($data.Length) data rows (rows $firstDataRow-$lastDataRow)."
